# This workbook ("paises.xlsx") tracks worldwide COVID-19 case counts, one
# country per row, sorted by total cases (column B) descending. This edit
# refreshes the data to the next scrape (23 May 2020 23:35 -> 24 May 2020
# 00:05), which updates case counts for several countries and - because the
# table stays sorted by total cases - causes a handful of neighboring rows
# to swap rank/order around the update boundary.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 00:05"

# --- Top of table: case-count refresh only (no re-ranking) ---------------
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1665674
$ws.Cells.Item(4, 3).Value = 20580
$ws.Cells.Item(4, 4).Value = 445286
$ws.Cells.Item(4, 5).Value = 1121739
$ws.Cells.Item(4, 7).Value = 1002
$ws.Cells.Item(4, 8).Value = 98649

# Row 16: Canada
$ws.Cells.Item(16, 2).Value = 83621
$ws.Cells.Item(16, 3).Value = 1141
$ws.Cells.Item(16, 4).Value = 43305
$ws.Cells.Item(16, 5).Value = 33961
$ws.Cells.Item(16, 7).Value = 105
$ws.Cells.Item(16, 8).Value = 6355

# --- Rows 126-130: re-ranked cluster around Sudan del Sur / Chad ---------
$ws.Cells.Item(126, 1).Value = "Sudan del Sur"
$ws.Cells.Item(126, 2).Value = 655
$ws.Cells.Item(126, 3).Value = 92
$ws.Cells.Item(126, 4).Value = 6
$ws.Cells.Item(126, 5).Value = 641
$ws.Cells.Item(126, 8).Value = 8

$ws.Cells.Item(127, 1).Value = "Republica del Chad"
$ws.Cells.Item(127, 2).Value = 648
$ws.Cells.Item(127, 3).Value = 37
$ws.Cells.Item(127, 4).Value = 204
$ws.Cells.Item(127, 5).Value = 384
$ws.Cells.Item(127, 7).Value = 2
$ws.Cells.Item(127, 8).Value = 60

$ws.Cells.Item(128, 1).Value = "Sierra Leona"
$ws.Cells.Item(128, 2).Value = 621
$ws.Cells.Item(128, 3).Value = 15
$ws.Cells.Item(128, 4).Value = 241
$ws.Cells.Item(128, 5).Value = 341
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = 39

$ws.Cells.Item(129, 1).Value = "Malta"
$ws.Cells.Item(129, 2).Value = 609
$ws.Cells.Item(129, 3).Value = 9
$ws.Cells.Item(129, 4).Value = 473
$ws.Cells.Item(129, 5).Value = 130
$ws.Cells.Item(129, 8).Value = 6

$ws.Cells.Item(130, 1).Value = "Nepal"
$ws.Cells.Item(130, 2).Value = 584
$ws.Cells.Item(130, 3).Value = 68
$ws.Cells.Item(130, 4).Value = 70
$ws.Cells.Item(130, 5).Value = 511
$ws.Cells.Item(130, 8).Value = 3

# --- Rows 147-148: Guayana Francesa / Nicaragua swap ----------------------
$ws.Cells.Item(147, 1).Value = "Guayana Francesa"
$ws.Cells.Item(147, 3).Value = 18
$ws.Cells.Item(147, 4).Value = 143
$ws.Cells.Item(147, 5).Value = 135
$ws.Cells.Item(147, 8).Value = 1

$ws.Cells.Item(148, 1).Value = "Nicaragua"
$ws.Cells.Item(148, 2).Value = 279
$ws.Cells.Item(148, 4).Value = 199
$ws.Cells.Item(148, 5).Value = 63
$ws.Cells.Item(148, 8).Value = 17

# --- Rows 155-157: Uganda / Martinica / Islas Feroe re-rank ---------------
$ws.Cells.Item(155, 1).Value = "Uganda"
$ws.Cells.Item(155, 2).Value = 198
$ws.Cells.Item(155, 3).Value = 23
$ws.Cells.Item(155, 4).Value = 68
$ws.Cells.Item(155, 5).Value = 130
$ws.Cells.Item(155, 8).Value = 0

$ws.Cells.Item(156, 1).Value = "Martinica"
$ws.Cells.Item(156, 2).Value = 197
$ws.Cells.Item(156, 4).Value = 91
$ws.Cells.Item(156, 5).Value = 92
$ws.Cells.Item(156, 8).Value = 14

$ws.Cells.Item(157, 1).Value = "Islas Feroe"
$ws.Cells.Item(157, 2).Value = 187
$ws.Cells.Item(157, 4).Value = 187
$ws.Cells.Item(157, 5).Value = 0

# --- Rows 199-200: Nueva Caledonia / Santa Lucia swap (name only) --------
$ws.Cells.Item(199, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(200, 1).Value = "Santa Lucia"

# --- Rows 209-211: Seychelles / Montserrat / Groenlandia re-rank ---------
$ws.Cells.Item(209, 1).Value = "Seychelles"

$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# --- Rows 214-216: Sahara Occidental / Bonaire... / San Bartolome swap ---
$ws.Cells.Item(214, 1).Value = "Sahara Occidental"
$ws.Cells.Item(215, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(216, 1).Value = "San Bartolome"
